$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Row 1: headers
$ws2.Range("A1").Value = "pseudo_member_id"
$ws2.Range("B1").Value = "session_id_mask"
$ws2.Range("C1").Value = "height_inches"
$ws2.Range("D1").Value = "weight_lbs"
$ws2.Range("E1").Value = "bp_systolic"
$ws2.Range("F1").Value = "bmi"

# Row 2
$ws2.Range("A2").Value = "dasda"
$ws2.Range("B2").Value = "asdas"
$ws2.Range("C2").Value = 128
$ws2.Range("D2").Value = 60
$ws2.Range("E2").Value = 140
$ws2.Range("F2").Value = 22

# Row 3
$ws2.Range("A3").Value = "asdasd"
$ws2.Range("B3").Value = "swdfsf"
$ws2.Range("C3").Value = 127
$ws2.Range("D3").Value = "NA"
$ws2.Range("E3").Value = "NA"
$ws2.Range("F3").Value = "NA"

# Row 12: headers (repeated)
$ws2.Range("A12").Value = "pseudo_member_id"
$ws2.Range("B12").Value = "session_id_mask"
$ws2.Range("C12").Value = "height_inches"
$ws2.Range("D12").Value = "weight_lbs"
$ws2.Range("E12").Value = "bp_systolic"
$ws2.Range("F12").Value = "bmi"

# Row 13
$ws2.Range("A13").Value = "dasda"
$ws2.Range("B13").Value = "asdas"
$ws2.Range("C13").Value = 1
$ws2.Range("D13").Value = 1
$ws2.Range("E13").Value = 1
$ws2.Range("F13").Value = 1

# Row 14
$ws2.Range("A14").Value = "asdasd"
$ws2.Range("B14").Value = "swdfsf"
$ws2.Range("C14").Value = 1
$ws2.Range("D14").Value = 0
$ws2.Range("E14").Value = 0
$ws2.Range("F14").Value = 0

# Column width for column A (target 28.81640625 chars; engine quantizes to
# the nearest 1/6-char pixel grid, so 28 is the input that lands closest)
$ws2.Columns.Item(1).ColumnWidth = 28

# Selection on sheet2
[void]$ws2.Range("C15").Select()
